$wb = $excel.ActiveWorkbook

# -- Overview sheet: row for 9ed7af9a-...md file (row 3) moves from
#    "Handed back: in sync with en-US" to "Ready for handoff"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-41-18 16:41:09"

# -- zh-cn sheet: row for 9ed7af9a-...md file (row 3) Status + Latest Handoff Datetime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-18 16:41:05"

# -- de-de sheet: row for 9ed7af9a-...md file (row 3) Status + Latest Handoff Datetime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-18 16:41:09"
